$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (Excel shifts existing rows 3:76 down to 4:77,
# and copies formatting from the row above by default).
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new record's values.
$ws.Cells.Item(3, 1).Value = 5
$ws.Cells.Item(3, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(3, 3).Value = "Maule"
$ws.Cells.Item(3, 4).Value = 44860
$ws.Cells.Item(3, 5).Value = 7
$ws.Cells.Item(3, 6).Value = 100112040
$ws.Cells.Item(3, 7).Value = "Cilantro"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 150
$ws.Cells.Item(3, 11).Value = 7000
$ws.Cells.Item(3, 12).Value = 7000
$ws.Cells.Item(3, 13).Value = 7000
$ws.Cells.Item(3, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(3, 15).Value = "Región del Maule"
$ws.Cells.Item(3, 16).Value = 194
$ws.Cells.Item(3, 17).Value = 36
$ws.Cells.Item(3, 18).Value = "Hortaliza"

Write-Host "Done"
